{"js": "const doc = context.document;\nconst body = doc.body;\n\n// 1. Remove the author/student-ID paragraph that follows the title\n//    (\"\u4f55\u6c83\u6d32 \u5b66\u53f72017210719 \"). Deleting the paragraph (including its\n//    paragraph mark) merges it away cleanly.\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nparas.items[1].delete();\nawait context.sync();\n\n// 2. Word re-stamps the \"_GoBack\" bookmark (last-edit location) at the\n//    point of the most recent edit whenever the file is saved. After the\n//    deletion above, that is the very start of the following heading\n//    paragraph (\"1 \u5b9e\u9a8c\u5185\u5bb9\"). Office.js does not auto-replace a same-named\n//    bookmark the way the Word COM object model does, so remove the old\n//    one explicitly before inserting the new one.\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst paras2 = body.paragraphs;\nparas2.load(\"items\");\nawait context.sync();\nconst headingRange = paras2.items[1].getRange(\"Start\");\nheadingRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3. Collapse the three separate runs that spell out the quotation\n//    (\"\u201c\", \"k rules \", \"\u201d\") into a single run \"\u201ck rules \u201d\" by searching\n//    for the full quoted phrase and replacing it with itself, which\n//    makes Word rebuild the matched span as one run.\nconst results = body.search(\"\\u201ck rules \\u201d\", { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\nresults.items[0].insertText(\"\\u201ck rules \\u201d\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the author/student-ID paragraph that follows the title\n#    (\"\u4f55\u6c83\u6d32 \u5b66\u53f72017210719 \"). Deleting the whole paragraph range\n#    (including its paragraph mark) merges it away cleanly, the same\n#    way Word does when you select the paragraph and hit Delete.\n$d.Paragraphs(2).Range.Delete()\n\n# 2. Word re-stamps the \"_GoBack\" bookmark (last-edit location) at the\n#    point of the most recent edit whenever the file is saved. After the\n#    deletion above, that is the very start of the following heading\n#    paragraph (\"1 \u5b9e\u9a8c\u5185\u5bb9\"). Adding a bookmark with the reserved name\n#    \"_GoBack\" automatically removes any pre-existing bookmark of that\n#    name elsewhere in the document (its old spot inside \"...\u6309\u4f4d|\u5f02\u6216...\").\n$target = $d.Paragraphs(2).Range.Duplicate\n$target.Collapse(1)\n$d.Bookmarks.Add(\"_GoBack\", $target)\n\n# 3. Collapse the three separate runs that spell out the quotation\n#    (\"\u201c\", \"k rules \", \"\u201d\") into a single run \"\u201ck rules \u201d\". Using\n#    Find/Replace (instead of touching Range.Text directly) makes Word\n#    rebuild the matched span as one run, mirroring the authored edit.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"\u201ck rules \u201d\"\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Replacement.Text = \"\u201ck rules \u201d\"\n$rng.Find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n"}
